$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a third header row (row 3) with French descriptions/labels for each
# of the existing columns (Operator, SampleID, Date, LaboratoryOperatingMode,
# CriticalApparatusCriticalSoftware, CriticalProduct, RawDataPathway,
# DilutionFactor, LivingCellCount).
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"

# The last two columns (DilutionFactor, LivingCellCount) get an explicit
# but blank description cell. A plain empty-string assignment clears the
# cell entirely, so force the cell to materialize as an (empty) text value
# via the quote-prefix, then strip the resulting formatting back off so no
# visible change/style is left behind - leaving a real, present, empty
# text cell like the template expects.
$ws.Range("H3").Value = "'"
$ws.Range("I3").Value = "'"
$ws.Range("H3:I3").ClearFormats()
